# Update loading_percent values for Case_2_55 (380 kV case)
# Applies the new computed loading percentages to the data rows (rows 2-25)
# across the affected columns (B, C, E, F, G, J, K, L, N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 21.33754624274098
    "C2" = 4.135940079205809
    "E2" = 11.01900336962196
    "F2" = 51.32974407439013
    "G2" = 3.772469011537285
    "J2" = 10.87603742891439
    "K2" = 17.34765440042642
    "L2" = 11.11243425105518
    "N2" = 25.19527429051162
    "B3" = 21.20996803257632
    "C3" = 3.941882621224618
    "E3" = 11.03494984647233
    "F3" = 51.24501399964127
    "G3" = 3.77564093418059
    "J3" = 10.89233416895806
    "K3" = 17.2638771304092
    "L3" = 11.12036145659792
    "N3" = 25.23264789273557
    "B4" = 21.13666461007873
    "C4" = 3.816859505132991
    "E4" = 11.04623756723457
    "F4" = 51.20274481621509
    "G4" = 3.77769015696669
    "J4" = 10.90282760609043
    "K4" = 17.21645631513574
    "L4" = 11.12708107015416
    "N4" = 25.25745003003835
    "B5" = 21.10808342398274
    "C5" = 3.764464734514821
    "E5" = 11.05121433979316
    "F5" = 51.187981338215
    "G5" = 3.778550883396352
    "J5" = 10.90722664466
    "K5" = 17.19815753396982
    "L5" = 11.13028573362685
    "N5" = 25.26802345268013
    "B6" = 21.10341619508105
    "C6" = 3.755678203940244
    "E6" = 11.05206351250563
    "F6" = 51.18567875708295
    "G6" = 3.778695358324449
    "J6" = 10.90796453417602
    "K6" = 17.19518138848375
    "L6" = 11.13084604791497
    "N6" = 25.26980732789404
    "B7" = 21.13627389648268
    "C7" = 3.816158706236629
    "E7" = 11.04630315872426
    "F7" = 51.20253573434974
    "G7" = 3.777701661041868
    "J7" = 10.90288643505066
    "K7" = 17.21620536008622
    "L7" = 11.12712240036827
    "N7" = 25.25759073873737
    "B8" = 21.29252932430287
    "C8" = 4.070263704082012
    "E8" = 11.02419149602352
    "F8" = 51.29851002788836
    "G8" = 3.773541647951203
    "J8" = 10.88155565820482
    "K8" = 17.3179433635649
    "L8" = 11.11478350371579
    "N8" = 25.2077757853756
    "B9" = 21.63763174447571
    "C9" = 4.520933607138216
    "E9" = 10.99267661223169
    "F9" = 51.56367372614571
    "G9" = 3.766186269765007
    "J9" = 10.84357429968304
    "K9" = 17.54861107270462
    "L9" = 11.10525380997054
    "N9" = 25.12480459620993
    "B10" = 21.91299261086929
    "C10" = 4.821904272693271
    "E10" = 10.97670510972501
    "F10" = 51.80469982671394
    "G10" = 3.761265629241763
    "J10" = 10.81799093168531
    "K10" = 17.73601347775944
    "L10" = 11.10714725440744
    "N10" = 25.07281780591999
    "B11" = 22.04258711002054
    "C11" = 4.952117625723457
    "E11" = 10.9709897314059
    "F11" = 51.92419449931908
    "G11" = 3.759130820672751
    "J11" = 10.80685128033625
    "K11" = 17.82490860564288
    "L11" = 11.10992750522756
    "N11" = 25.05111631922987
    "B12" = 22.09224743339025
    "C12" = 5.000452573158848
    "E12" = 10.96904752451119
    "F12" = 51.97084142965353
    "G12" = 3.75833722940887
    "J12" = 10.80270425783332
    "K12" = 17.85907206466082
    "L12" = 11.11125492636194
    "N12" = 25.04317867009377
    "B13" = 22.08152674793774
    "C13" = 4.990086241812653
    "E13" = 10.96945594997319
    "F13" = 51.96073338287896
    "G13" = 3.758507485966422
    "J13" = 10.8035942266668
    "K13" = 17.85169246017912
    "L13" = 11.11095685236179
    "N13" = 25.04487571854188
    "B14" = 22.04666117303923
    "C14" = 4.9561137443633
    "E14" = 10.97082549889903
    "F14" = 51.92800429332713
    "G14" = 3.759065235026002
    "J14" = 10.80650867484735
    "K14" = 17.82770935397652
    "L14" = 11.1100312177294
    "N14" = 25.05045766674772
    "B15" = 22.0253801559202
    "C15" = 4.935177486986307
    "E15" = 10.97169328360424
    "F15" = 51.90813808445925
    "G15" = 3.75940879917819
    "J15" = 10.8083031366416
    "K15" = 17.81308352040837
    "L15" = 11.10949995874477
    "N15" = 25.05391327111043
    "B16" = 21.90460731423576
    "C16" = 4.813258747022449
    "E16" = 10.97710974643246
    "F16" = 51.79708715318105
    "G16" = 3.761407221702185
    "J16" = 10.81872893301622
    "K16" = 17.73027534364638
    "L16" = 11.10700406092771
    "N16" = 25.07427525393007
    "B17" = 21.83160039685891
    "C17" = 4.736741274382414
    "E17" = 10.98082902645934
    "F17" = 51.73147060968649
    "G17" = 3.762659666701053
    "J17" = 10.82525221779918
    "K17" = 17.68039245182382
    "L17" = 11.10596348270745
    "N17" = 25.08726561148275
    "B18" = 21.79001874769317
    "C18" = 4.692100309014537
    "E18" = 10.9831142166836
    "F18" = 51.69465819353523
    "G18" = 3.763389797793988
    "J18" = 10.82905116602682
    "K18" = 17.65204566075337
    "L18" = 11.10554567600957
    "N18" = 25.09492062472861
    "B19" = 21.77601144737534
    "C19" = 4.676877780859476
    "E19" = 10.98391303573916
    "F19" = 51.68235417842804
    "G19" = 3.763638686165531
    "J19" = 10.83034549485624
    "K19" = 17.64250779241519
    "L19" = 11.10543528869689
    "N19" = 25.09754395790844
    "B20" = 21.83932994034169
    "C20" = 4.744951991214823
    "E20" = 10.9804180024007
    "F20" = 51.738359644925
    "G20" = 3.762525332514861
    "J20" = 10.82455294869664
    "K20" = 17.68566707939424
    "L20" = 11.10605556085496
    "N20" = 25.08586379342918
    "B21" = 22.05688644322506
    "C21" = 4.966118798167145
    "E21" = 10.97041720872963
    "F21" = 51.93757987430535
    "G21" = 3.758901009226349
    "J21" = 10.8056506978182
    "K21" = 17.83474037237434
    "L21" = 11.11029565784916
    "N21" = 25.04881050637525
    "B22" = 22.20246815975867
    "C22" = 5.104983382634728
    "E22" = 10.96517516416158
    "F22" = 52.07591361367714
    "G22" = 3.756618614893644
    "J22" = 10.79371254563097
    "K22" = 17.93507601282671
    "L22" = 11.11466666630915
    "N22" = 25.02622753940409
    "B23" = 22.12447026611363
    "C23" = 5.031391447478287
    "E23" = 10.96785482103481
    "F23" = 52.0013451742871
    "G23" = 3.757828902373415
    "J23" = 10.80004625002592
    "K23" = 17.88126688095601
    "L23" = 11.11218786841225
    "N23" = 25.03813096902917
    "B24" = 21.83583419300161
    "C24" = 4.741241950116279
    "E24" = 10.98060336867552
    "F24" = 51.73524227165294
    "G24" = 3.762586033583638
    "J24" = 10.82486893676334
    "K24" = 17.68328138577385
    "L24" = 11.10601337020159
    "N24" = 25.08649697384813
    "B25" = 21.54030666196599
    "C25" = 4.404242486107412
    "E25" = 10.99993789815655
    "F25" = 51.48376634615395
    "G25" = 3.768090794540081
    "J25" = 10.85343982185936
    "K25" = 17.4829783191409
    "L25" = 11.10626590330503
    "N25" = 25.14567500005326
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
